$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 192.25
$ws.Range("I33").Value = 170.8
$ws.Range("J33").Value = 299.5
$ws.Range("K33").Value = 170.8
$ws.Range("L33").Value = 299.5
$ws.Range("M33").Value = 58.19999999999999
$ws.Range("N33").Value = -757.5
$ws.Range("H86").Value = 5801
$ws.Range("I86").Value = 5701.5
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 5701.5
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -4578.5
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5801
$ws.Range("I89").Value = 5701.5
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 28507.5
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -22891.5
$ws.Range("N89").Value = -41232
$ws.Range("H129").Value = 2589.4285
$ws.Range("J129").Value = 4523.4287
$ws.Range("L129").Value = 13570.2861
$ws.Range("N129").Value = -23570.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1252.4445
$ws.Range("I74").Value = 858.5
$ws.Range("J74").Value = 4404
$ws.Range("K74").Value = 858.5
$ws.Range("L74").Value = 4404
$ws.Range("M74").Value = 15.5
$ws.Range("N74").Value = -6152
$ws.Range("H77").Value = 1252.4445
$ws.Range("I77").Value = 858.5
$ws.Range("J77").Value = 4404
$ws.Range("K77").Value = 4292.5
$ws.Range("L77").Value = 22020
$ws.Range("M77").Value = 75.5
$ws.Range("N77").Value = -30756
$ws.Range("H122").Value = 315304.66
$ws.Range("I122").Value = 437561.1
$ws.Range("K122").Value = 1312683.3
$ws.Range("M122").Value = -1310233.3
$ws.Range("H132").Value = 1536.2916
$ws.Range("I132").Value = 1434.3334
$ws.Range("K132").Value = 4303.0002
$ws.Range("M132").Value = -1773.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6849.75
$ws.Range("I20").Value = 7219.7
$ws.Range("K20").Value = 7219.7
$ws.Range("M20").Value = -6972.7
$ws.Range("H94").Value = 769.55
$ws.Range("I94").Value = 743.94446
$ws.Range("K94").Value = 743.94446
$ws.Range("M94").Value = -292.94446
$ws.Range("H107").Value = 1278
$ws.Range("I107").Value = 1330
$ws.Range("K107").Value = 1330
$ws.Range("M107").Value = 590
$ws.Range("H134").Value = 2482.9395
$ws.Range("I134").Value = 2298.182
$ws.Range("K134").Value = 6894.545999999999
$ws.Range("M134").Value = -4359.545999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 45865
$ws.Range("J74").Value = 44988
$ws.Range("L74").Value = 44988
$ws.Range("N74").Value = -46736
$ws.Range("H77").Value = 45865
$ws.Range("J77").Value = 44988
$ws.Range("L77").Value = 134964
$ws.Range("N77").Value = -143700
$ws.Range("H99").Value = 12827.167
$ws.Range("I99").Value = 8760.637000000001
$ws.Range("J99").Value = 16268.077
$ws.Range("K99").Value = 8760.637000000001
$ws.Range("L99").Value = 16268.077
$ws.Range("M99").Value = -7262.637000000001
$ws.Range("N99").Value = -19264.077
$ws.Range("H103").Value = 17678
$ws.Range("I103").Value = 17678
$ws.Range("K103").Value = 17678
$ws.Range("M103").Value = -16506
$ws.Range("H126").Value = 12827.167
$ws.Range("I126").Value = 8760.637000000001
$ws.Range("J126").Value = 16268.077
$ws.Range("K126").Value = 26281.911
$ws.Range("L126").Value = 48804.231
$ws.Range("M126").Value = -23811.911
$ws.Range("N126").Value = -53744.231
$ws.Range("H134").Value = 1965.12
$ws.Range("I134").Value = 1549.8918
$ws.Range("J134").Value = 3146.923
$ws.Range("K134").Value = 4649.6754
$ws.Range("L134").Value = 9440.769
$ws.Range("M134").Value = -2114.6754
$ws.Range("N134").Value = -14510.769
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 997.4706
$ws.Range("I14").Value = 997.4706
$ws.Range("K14").Value = 2992.4118
$ws.Range("M14").Value = -2819.4118
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H86").Value = 345
$ws.Range("I86").Value = 345
$ws.Range("K86").Value = 1035
$ws.Range("M86").Value = 151
$ws.Range("H89").Value = 345
$ws.Range("I89").Value = 345
$ws.Range("K89").Value = 3105
$ws.Range("M89").Value = 2823
$ws.Range("H92").Value = 444.6842
$ws.Range("I92").Value = 483.75
$ws.Range("J92").Value = 416.27274
$ws.Range("K92").Value = 1451.25
$ws.Range("L92").Value = 1248.81822
$ws.Range("M92").Value = -203.25
$ws.Range("N92").Value = -3744.81822
$ws.Range("H95").Value = 2350
$ws.Range("J95").Value = 3200
$ws.Range("L95").Value = 9600
$ws.Range("N95").Value = -13718
$ws.Range("H107").Value = 383.25
$ws.Range("J107").Value = 427.625
$ws.Range("L107").Value = 1282.875
$ws.Range("N107").Value = -5122.875
$ws.Range("H119").Value = 4042.25
$ws.Range("I119").Value = 3066.6667
$ws.Range("K119").Value = 9200.000100000001
$ws.Range("M119").Value = -4362.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7268.091
$ws.Range("I70").Value = 7108.1665
$ws.Range("J70").Value = 7460
$ws.Range("K70").Value = 7108.1665
$ws.Range("L70").Value = 7460
$ws.Range("M70").Value = -6838.1665
$ws.Range("N70").Value = -8000
$ws.Range("H73").Value = 7268.091
$ws.Range("I73").Value = 7108.1665
$ws.Range("J73").Value = 7460
$ws.Range("K73").Value = 7108.1665
$ws.Range("L73").Value = 7460
$ws.Range("M73").Value = -6172.1665
$ws.Range("N73").Value = -9332
$ws.Range("H80").Value = 3823.65
$ws.Range("I80").Value = 3177.2144
$ws.Range("J80").Value = 5332
$ws.Range("K80").Value = 3177.2144
$ws.Range("L80").Value = 5332
$ws.Range("M80").Value = -2179.2144
$ws.Range("N80").Value = -7328
$ws.Range("H83").Value = 3823.65
$ws.Range("I83").Value = 3177.2144
$ws.Range("J83").Value = 5332
$ws.Range("K83").Value = 15886.072
$ws.Range("L83").Value = 26660
$ws.Range("M83").Value = -10894.072
$ws.Range("N83").Value = -36644
$ws.Range("H113").Value = 4436
$ws.Range("I113").Value = 6094
$ws.Range("J113").Value = 2778
$ws.Range("K113").Value = 6094
$ws.Range("L113").Value = 2778
$ws.Range("M113").Value = -3924
$ws.Range("N113").Value = -7118
$ws.Range("H126").Value = 3760
$ws.Range("I126").Value = 2813.4
$ws.Range("J126").Value = 5337.6665
$ws.Range("K126").Value = 8440.200000000001
$ws.Range("L126").Value = 16012.9995
$ws.Range("M126").Value = -5970.200000000001
$ws.Range("N126").Value = -20952.9995
$ws.Range("H132").Value = 2191.7942
$ws.Range("I132").Value = 1711.4762
$ws.Range("J132").Value = 2967.6924
$ws.Range("K132").Value = 5134.4286
$ws.Range("L132").Value = 8903.0772
$ws.Range("M132").Value = -2604.4286
$ws.Range("N132").Value = -13963.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4424.3335
$ws.Range("I7").Value = 4102.375
$ws.Range("K7").Value = 4102.375
$ws.Range("M7").Value = -3990.375
$ws.Range("H22").Value = 25299.5
$ws.Range("J22").Value = 49999
$ws.Range("L22").Value = 49999
$ws.Range("N22").Value = -50589
$ws.Range("H27").Value = 25299.5
$ws.Range("J27").Value = 49999
$ws.Range("L27").Value = 49999
$ws.Range("N27").Value = -50213
$ws.Range("H40").Value = 3098.6
$ws.Range("I40").Value = 3098.6
$ws.Range("K40").Value = 3098.6
$ws.Range("M40").Value = -2962.6
$ws.Range("H61").Value = 5999.6
$ws.Range("J61").Value = 3000
$ws.Range("L61").Value = 3000
$ws.Range("N61").Value = -3404
$ws.Range("H113").Value = 5999.6
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340
$ws.Range("H122").Value = 7513.5
$ws.Range("J122").Value = 5030.4
$ws.Range("L122").Value = 15091.2
$ws.Range("N122").Value = -19991.2
$ws.Range("H126").Value = 4424.3335
$ws.Range("I126").Value = 4102.375
$ws.Range("K126").Value = 12307.125
$ws.Range("M126").Value = -9837.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8165.643
$ws.Range("I62").Value = 6500
$ws.Range("J62").Value = 8443.25
$ws.Range("K62").Value = 6500
$ws.Range("L62").Value = 8443.25
$ws.Range("M62").Value = -5876
$ws.Range("N62").Value = -9691.25
$ws.Range("H65").Value = 8165.643
$ws.Range("I65").Value = 6500
$ws.Range("J65").Value = 8443.25
$ws.Range("K65").Value = 32500
$ws.Range("L65").Value = 42216.25
$ws.Range("M65").Value = -29380
$ws.Range("N65").Value = -48456.25
$ws.Range("H74").Value = 60899.332
$ws.Range("J74").Value = 60899.332
$ws.Range("L74").Value = 60899.332
$ws.Range("N74").Value = -62771.332
$ws.Range("H77").Value = 60899.332
$ws.Range("J77").Value = 60899.332
$ws.Range("L77").Value = 182697.996
$ws.Range("N77").Value = -192057.996
$ws.Range("H126").Value = 2145.4707
$ws.Range("I126").Value = 1654.5625
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 4963.6875
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -2493.6875
$ws.Range("N126").Value = -34940
$ws.Range("H136").Value = 1339.64
$ws.Range("I136").Value = 1325.7391
$ws.Range("J136").Value = 1499.5
$ws.Range("K136").Value = 3977.2173
$ws.Range("L136").Value = 4498.5
$ws.Range("M136").Value = -1427.2173
$ws.Range("N136").Value = -9598.5
